$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table S1
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table S1")
$ws1.Range("B4").Value  = "SMWD < ref., %"
$ws1.Range("B13").Value = "Imp. general health (VAS < 73, EQ5D5L), %"
$ws1.Range("B58").Value = "reduced LVEF, %"

# ---------------------------------------------------------------------
# Table S4
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table S4")
$ws4.Range("A15").Value = "CT abnormality (CT score " + [char]8805 + " 1)"

# ---------------------------------------------------------------------
# Table S5
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Table S5")
$ws5.Range("A5").Value  = "SMWD < ref."
$ws5.Range("A6").Value  = "Fatigue score (likert CFS)"
$ws5.Range("A8").Value  = "General health score (EQ5D5L VAS)"
$ws5.Range("A9").Value  = "Imp. general health (VAS < 73, EQ5D5L)"
$ws5.Range("A10").Value = "Mobility impairment score (EQ5D5L)"
$ws5.Range("A11").Value = "Imp. mobility (score  > 1, EQ5D5L)"
$ws5.Range("A12").Value = "Self-care impairment score (EQ5D5L)"
$ws5.Range("A13").Value = "Imp. self-care (score  > 1, EQ5D5L)"
$ws5.Range("A14").Value = "Activity impairment score (EQ5D5L)"
$ws5.Range("A15").Value = "Imp. usual activity (score  > 1, EQ5D5L)"
$ws5.Range("A16").Value = "Pain/discomfort score (EQ5D5L)"
$ws5.Range("A17").Value = "Pain/discomfort (score  > 1, EQ5D5L)"
$ws5.Range("A18").Value = "Anxiety/depression score (EQ5D5L)"
$ws5.Range("A19").Value = "Anxiety/depression (score  > 1, EQ5D5L)"
$ws5.Range("A20").Value = "Stress score (PSS)"
$ws5.Range("A22").Value = "Somatic symptom disorder score (SSD-12)"
$ws5.Range("A23").Value = "Resilience score (BRCS)"

$ws5.Range("B22").Value = "median: 4 [IQR: 1 - 7.5]" + [char]10 + "range: 0 - 24"
$ws5.Range("C22").Value = "median: 5 [IQR: 2 - 10]" + [char]10 + "range: 0 - 26"
$ws5.Range("D22").Value = "median: 20 [IQR: 12 - 24]" + [char]10 + "range: 7 - 30"

$ws5.Range("B23").Value = "median: 16 [IQR: 12 - 19]" + [char]10 + "range: 4 - 20"
$ws5.Range("C23").Value = "median: 16 [IQR: 12 - 18]" + [char]10 + "range: 4 - 20"
$ws5.Range("D23").Value = "median: 16 [IQR: 14 - 18]" + [char]10 + "range: 9 - 19"

$ws5.Range("B24").Value = "low: 35% (n = 12)" + [char]10 + "medium: 24% (n = 8)" + [char]10 + "high: 41% (n = 14)" + [char]10 + "n = 34"
$ws5.Range("C24").Value = "low: 31% (n = 10)" + [char]10 + "medium: 19% (n = 6)" + [char]10 + "high: 50% (n = 16)" + [char]10 + "n = 32"
$ws5.Range("D24").Value = "low: 17% (n = 3)" + [char]10 + "medium: 39% (n = 7)" + [char]10 + "high: 44% (n = 8)" + [char]10 + "n = 18"
